$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers: OutSample -> TestSample
$ws.Range("E1").Value = "RMSE_TestSample"
$ws.Range("F1").Value = "R2_TestSample"
$ws.Range("G1").Value = "Adjusted_R2_TestSample"

# Updated results table (rows 2-9, columns B-J)
$data = @(
    @(0.04967061671052313, 0.8024266576305294, 0.7982050050157972, 0.07367127193719293, 0.9637309116279517, 0.9603726627046139, 131.0858086492545, -252.1716172985091, -234.7684226817991),
    @(0.0430195623246154, 0.8322296691693302, 0.8279094031393559, 0.06969859823890841, 0.9665727112180673, 0.9627884898465278, 161.7198899650198, -311.4397799300395, -290.5559463899876),
    @(0.04246927223592505, 0.847015936123746, 0.8424000376447212, 0.07005820289497666, 0.9667131045806088, 0.9622321763510754, 164.7879603469766, -315.5759206939533, -291.2114482305593),
    @(0.0422367466578034, 0.8541423743941738, 0.8490910280528465, 0.0697987185196205, 0.9669380016364926, 0.9617518058147659, 166.5517374879095, -317.103474975819, -289.2583635890831),
    @(0.04037852781878488, 0.8716163745054414, 0.8665926674208717, 0.06820949929607642, 0.9684479120228996, 0.9627685361870215, 178.0131686143943, -338.0263372287885, -306.7005869187106),
    @(0.03778366543972687, 0.8826457640777204, 0.8775211249544768, 0.064231858796339, 0.9710139330818076, 0.965098409220952, 186.4478464623334, -352.8956929246668, -318.0893036912469),
    @(0.03655437199078362, 0.8910743815879055, 0.8858191982434623, 0.06363285793608846, 0.9713504582728413, 0.9647849382937007, 191.0822397680694, -360.1644795361389, -321.877451379377),
    @(0.03602200881466022, 0.895574745970436, 0.8900544682243797, 0.06254077995352933, 0.9741104121661414, 0.9675003046340923, 194.0576474689833, -364.1152949379666, -322.3476278578627)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = 2 + $j  # Column B = 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

$wb.Save()
